# "Removed parts related to workflow statistics."
#
# The underlying sheet1 row 12 ("Undocumented unreachable activities") had
# its Action column mis-set to "Double check" -- correct it to "Fix".
# The workbook also had the "Project" sheet as the active/selected tab with
# a stale selection; re-point the active tab back to "Workflow" (clearing
# its old selection) and update the "Project" sheet's leftover selection.

$wb = $excel.ActiveWorkbook

$wsWorkflow = $wb.Worksheets.Item("Workflow")
$wsProject  = $wb.Worksheets.Item("Project")

# Fix the mis-categorized check in the Workflow sheet.
$wsWorkflow.Range("E12").Value = "Fix"

# Update the (now inactive) Project sheet's remembered selection.
$wsProject.Activate() | Out-Null
$wsProject.Range("C6").Select() | Out-Null

# Make Workflow the active sheet/tab again, with a fresh selection.
$wsWorkflow.Activate() | Out-Null
$wsWorkflow.Range("A1").Select() | Out-Null
